# ============================================================================
# Test_Suite_Execution_Report_Analysis.xlsx edit script
# Commit: "Features added - Execution on Safari Browser, ExtentReport name"
# ============================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename the worksheet
# ---------------------------------------------------------------------------
$ws.Name = "Test_Execution"

# ---------------------------------------------------------------------------
# 2. Wipe the old layout completely (values + number formats + fonts +
#    fills + borders + alignment) so we can rebuild the new table from a
#    clean slate.
# ---------------------------------------------------------------------------
$ws.Range("A1:M20").Clear()

# Old column B held the long TestNG class names - no longer used, reset it
# back to the sheet's standard width.
$ws.Columns("B").ColumnWidth = $ws.StandardWidth

# ---------------------------------------------------------------------------
# 3. New column widths for the relocated table (now starting at column C)
# ---------------------------------------------------------------------------
$ws.Columns("C").ColumnWidth = 24.09
$ws.Columns("D").ColumnWidth = 12.27

# ---------------------------------------------------------------------------
# 4. Row heights
# ---------------------------------------------------------------------------
$ws.Rows(5).RowHeight = 15
$ws.Rows(6).RowHeight = 18.5
$ws.Rows(7).RowHeight = 29
$ws.Rows(8).RowHeight = 29.5

# ---------------------------------------------------------------------------
# 5. Table values
# ---------------------------------------------------------------------------
# Header row (row 6)
$ws.Range("C6").Value = "Mode"
$ws.Range("D6").Value = "OS"
$ws.Range("E6").Value = "Chrome"
$ws.Range("F6").Value = "Edge"
$ws.Range("G6").Value = "Opera"
$ws.Range("H6").Value = "Headless"
$ws.Range("I6").Value = "Firefox"
$ws.Range("J6").Value = "Safari"

# Sequential execution row (row 7)
$ws.Range("C7").Value = "Sequential" + [char]10 + "(1 browser - One after One)"
$ws.Range("D7").Value = "WINDOWS 10"
$ws.Range("E7").Value = "2 min 48 sec"
$ws.Range("F7").Value = "2 min 18 sec"
$ws.Range("G7").Value = "2 min 33 sec"
$ws.Range("H7").Value = "2 min 14 sec"

# Parallel execution row (row 8)
$ws.Range("C8").Value = "Parallel " + [char]10 + "(1 browser - 5 instances)"
$ws.Range("D8").Value = "WINDOWS 10"
$ws.Range("E8").Value = "1 min 14 sec"
$ws.Range("F8").Value = "0 min 56 sec"
$ws.Range("G8").Value = "1 min 4 sec"
$ws.Range("H8").Value = "1 min 18 sec"

# ---------------------------------------------------------------------------
# 6. Formatting - borders
#    Whole table: thin grid everywhere, then a medium (thick) box around the
#    outside - matches the look used on the previous table layout.
# ---------------------------------------------------------------------------
$table = $ws.Range("C6:J8")
$table.Borders.LineStyle = 1
$table.Borders.Weight = 2
$table.BorderAround(1, -4138)

# ---------------------------------------------------------------------------
# 7. Formatting - header row (bold, bigger font, yellow fill)
# ---------------------------------------------------------------------------
$header = $ws.Range("C6:J6")
$header.Font.Bold = $true
$header.Font.Size = 14
$header.Interior.Color = 65535

# ---------------------------------------------------------------------------
# 8. Formatting - wrap text on the mode labels and the OS column
# ---------------------------------------------------------------------------
$ws.Range("C7").WrapText = $true
$ws.Range("C8").WrapText = $true

$osCol = $ws.Range("D7:D8")
$osCol.WrapText = $true
$osCol.Font.Bold = $true

# ---------------------------------------------------------------------------
# 9. Selection / view state
# ---------------------------------------------------------------------------
$ws.Range("E12").Select()
